$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Parag%Goyal%NULL%0,           Justin J.%Choi%NULL%0,           Laura C.%Pinheiro%NULL%0,           Edward J.%Schenck%NULL%0,           Ruijun%Chen%NULL%0,           Assem%Jabri%NULL%0,           Michael J.%Satlin%NULL%0,           Thomas R.%Campion%NULL%0,           Musarrat%Nahid%NULL%0,           Joanna B.%Ringel%NULL%0,           Katherine L.%Hoffman%NULL%0,           Mark N.%Alshak%NULL%0,           Han A.%Li%NULL%0,           Graham T.%Wehmeyer%NULL%0,           Graham T.%Wehmeyer%NULL%0,           Mangala%Rajan%NULL%0,           Evgeniya%Reshetnyak%NULL%0,           Nathaniel%Hupert%NULL%0,           Evelyn M.%Horn%NULL%0,           Fernando J.%Martinez%NULL%0,           Roy M.%Gulick%NULL%0,           Monika M.%Safford%NULL%0,           Monika M.%Safford%NULL%0]'
$ws.Range('E3').Value = '[Chaolin%Huang%NULL%0,           Yeming%Wang%NULL%0,           Xingwang%Li%NULL%0,           Lili%Ren%NULL%0,           Jianping%Zhao%NULL%0,           Yi%Hu%NULL%0,           Li%Zhang%NULL%0,           Guohui%Fan%NULL%0,           Jiuyang%Xu%NULL%0,           Xiaoying%Gu%NULL%0,           Zhenshun%Cheng%NULL%0,           Ting%Yu%NULL%0,           Jiaan%Xia%NULL%0,           Yuan%Wei%NULL%0,           Wenjuan%Wu%NULL%0,           Xuelei%Xie%NULL%0,           Wen%Yin%NULL%0,           Hui%Li%NULL%0,           Min%Liu%NULL%0,           Yan%Xiao%NULL%0,           Hong%Gao%NULL%0,           Li%Guo%NULL%0,           Jungang%Xie%NULL%0,           Guangfa%Wang%NULL%0,           Rongmeng%Jiang%NULL%0,           Zhancheng%Gao%NULL%0,           Qi%Jin%NULL%0,           Jianwei%Wang%wangjw28@163.com%0,           Bin%Cao%caobin_ben@163.com%0]'
$ws.Range('E4').Value = '[Suxin%Wan%NULL%0,           Yi%Xiang%NULL%0,           Wei%Fang%NULL%0,           Yu%Zheng%yuzheng1@cdutcm.edu.cn%0,           Boqun%Li%279685211@qq.com%0,           Boqun%Li%279685211@qq.com%0,           Yanjun%Hu%huyanjun@163.com%0,           Chunhui%Lang%NULL%0,           Daoqiu%Huang%NULL%0,           Qiuyan%Sun%NULL%0,           Yan%Xiong%NULL%0,           Xia%Huang%NULL%0,           Jinglong%Lv%NULL%0,           Yaling%Luo%NULL%0,           Li%Shen%NULL%0,           Haoran%Yang%NULL%0,           Gu%Huang%NULL%0,           Ruishan%Yang%NULL%0]'
$ws.Range('E5').Value = '[Ruirui%Wang%NULL%0,           Min%Pan%NULL%0,           Xiumei%Zhang%NULL%0,           Mingfeng%Han%fyhmf@163.com%0,           Xiaoyun%Fan%13956988552@126.com%0,           Fengde%Zhao%NULL%0,           Manli%Miao%NULL%0,           Jing%Xu%NULL%0,           Minglong%Guan%NULL%0,           Xia%Deng%NULL%0,           Xu%Chen%NULL%0,           Leilei%Shen%NULL%0]'
$ws.Range('E6').Value = '[Xiaoli%Zhang%NULL%0,           Huan%Cai%NULL%0,           Jianhua%Hu%NULL%0,           Jiangshan%Lian%NULL%0,           Jueqing%Gu%NULL%0,           Shanyan%Zhang%NULL%0,           Chanyuan%Ye%NULL%1,           Yingfeng%Lu%NULL%0,           Ciliang%Jin%NULL%0,           Guodong%Yu%NULL%0,           Hongyu%Jia%NULL%0,           Yimin%Zhang%NULL%0,           Jifang%Sheng%jifang_sheng@zju.edu.cn%0,           Lanjuan%Li%ljli@zju.edu.cn%2,           Yida%Yang%yidayang65@zju.edu.cn%0]'
$ws.Range('E7').Value = '[Fei%Zhou%NULL%0,           Ting%Yu%NULL%0,           Ronghui%Du%NULL%0,           Guohui%Fan%NULL%0,           Ying%Liu%NULL%0,           Zhibo%Liu%NULL%0,           Jie%Xiang%NULL%0,           Yeming%Wang%NULL%0,           Bin%Song%NULL%0,           Xiaoying%Gu%NULL%0,           Lulu%Guan%NULL%0,           Yuan%Wei%NULL%0,           Hui%Li%NULL%0,           Xudong%Wu%NULL%0,           Jiuyang%Xu%NULL%0,           Shengjin%Tu%NULL%0,           Yi%Zhang%NULL%0,           Hua%Chen%NULL%0,           Bin%Cao%NULL%0]'
$ws.Range('E8').Value = '[NULL%NULL%NULL%0,           NULL%NULL%NULL%0,           NULL%NULL%NULL%0,           Nancy%Chow%NULL%0,           Katherine%Fleming-Dutra%NULL%0,           Ryan%Gierke%NULL%0,           Aron%Hall%NULL%0,           Michelle%Hughes%NULL%0,           Tamara%Pilishvili%NULL%0,           Matthew%Ritchey%NULL%0,           Katherine%Roguski%NULL%0,           Tami%Skoff%NULL%0,           Emily%Ussery%NULL%0]'
$ws.Range('E9').Value = '[Tao%Chen%NULL%0,           Di%Wu%NULL%0,           Huilong%Chen%NULL%0,           Weiming%Yan%NULL%0,           Danlei%Yang%NULL%0,           Guang%Chen%NULL%0,           Ke%Ma%NULL%0,           Dong%Xu%NULL%0,           Haijing%Yu%NULL%0,           Hongwu%Wang%NULL%0,           Tao%Wang%NULL%0,           Wei%Guo%NULL%0,           Jia%Chen%NULL%0,           Chen%Ding%NULL%0,           Xiaoping%Zhang%NULL%0,           Jiaquan%Huang%NULL%0,           Meifang%Han%NULL%0,           Shusheng%Li%NULL%0,           Xiaoping%Luo%NULL%0,           Jianping%Zhao%NULL%0,           Qin%Ning%NULL%0]'
$ws.Range('E10').Value = '[Wei-jie%Guan%NULL%0,           Zheng-yi%Ni%NULL%0,           Zheng-yi%Ni%NULL%0,           Yu%Hu%NULL%0,           Wen-hua%Liang%NULL%0,           Chun-quan%Ou%NULL%0,           Jian-xing%He%NULL%0,           Lei%Liu%NULL%0,           Hong%Shan%NULL%0,           Chun-liang%Lei%NULL%0,           David S.C.%Hui%NULL%0,           Bin%Du%NULL%0,           Lan-juan%Li%NULL%0,           Guang%Zeng%NULL%0,           Kwok-Yung%Yuen%NULL%0,           Ru-chong%Chen%NULL%0,           Chun-li%Tang%NULL%0,           Tao%Wang%NULL%0,           Ping-yan%Chen%NULL%0,           Jie%Xiang%NULL%0,           Shi-yue%Li%NULL%0,           Jin-lin%Wang%NULL%0,           Zi-jing%Liang%NULL%0,           Yi-xiang%Peng%NULL%0,           Li%Wei%NULL%0,           Yong%Liu%NULL%0,           Ya-hua%Hu%NULL%0,           Peng%Peng%NULL%0,           Jian-ming%Wang%NULL%0,           Ji-yang%Liu%NULL%0,           Zhong%Chen%NULL%0,           Gang%Li%NULL%0,           Zhi-jian%Zheng%NULL%0,           Shao-qin%Qiu%NULL%0,           Jie%Luo%NULL%0,           Chang-jiang%Ye%NULL%0,           Shao-yong%Zhu%NULL%0,           Nan-shan%Zhong%NULL%0]'
$ws.Range('E11').Value = '[Xiaochen%Li%NULL%0,           Shuyun%Xu%NULL%0,           Muqing%Yu%NULL%0,           Ke%Wang%NULL%0,           Yu%Tao%NULL%0,           Ying%Zhou%NULL%0,           Jing%Shi%NULL%0,           Min%Zhou%NULL%0,           Bo%Wu%NULL%0,           Zhenyu%Yang%NULL%0,           Cong%Zhang%NULL%0,           Junqing%Yue%NULL%0,           Zhiguo%Zhang%NULL%0,           Harald%Renz%NULL%0,           Xiansheng%Liu%NULL%0,           Jungang%Xie%NULL%0,           Min%Xie%NULL%0,           Jianping%Zhao%NULL%0]'
$ws.Range('E12').Value = '[Akdis%Cezmi A%coreGivesNoEmail%0,          Cao%Yi-yuan%coreGivesNoEmail%0,          Dong%Xiang%coreGivesNoEmail%0,          Gao%Ya-dong%coreGivesNoEmail%0,          Yan%You-qin%coreGivesNoEmail%0,          Yang%Yi-bin%coreGivesNoEmail%0,          Yuan%Ya-dong%coreGivesNoEmail%0,          Zhang%Jin-jin%coreGivesNoEmail%0]'
$ws.Range('E13').Value = '[Yun%Feng%NULL%0,           Yun%Ling%NULL%0,           Yun%Ling%NULL%0,           Tao%Bai%NULL%0,           Tao%Bai%NULL%0,           Yusang%Xie%NULL%0,           Jie%Huang%NULL%0,           Jie%Huang%NULL%0,           Jian%Li%NULL%0,           Weining%Xiong%NULL%0,           Dexiang%Yang%NULL%0,           Rong%Chen%NULL%0,           Fangying%Lu%NULL%0,           Yunfei%Lu%NULL%0,           Xuhui%Liu%NULL%0,           Yuqing%Chen%NULL%0,           Yuqing%Chen%NULL%0,           Xin%Li%NULL%0,           Yong%Li%NULL%0,           Hanssa Dwarka%Summah%NULL%0,           Huihuang%Lin%NULL%0,           Jiayang%Yan%NULL%0,           Min%Zhou%NULL%0,           Hongzhou%Lu%NULL%0,           Hongzhou%Lu%NULL%0,           Jieming%Qu%NULL%0,           Jieming%Qu%NULL%0]'
$ws.Range('E14').Value = '[Eu Suk%Kim%NULL%1,           Bum Sik%Chin%NULL%2,           Bum Sik%Chin%NULL%0,           Chang Kyung%Kang%NULL%2,           Chang Kyung%Kang%NULL%0,           Nam Joong%Kim%NULL%2,           Nam Joong%Kim%NULL%0,           Yu Min%Kang%NULL%2,           Yu Min%Kang%NULL%0,           Jae-Phil%Choi%NULL%2,           Jae-Phil%Choi%NULL%0,           Dong Hyun%Oh%NULL%2,           Dong Hyun%Oh%NULL%0,           Jeong-Han%Kim%NULL%2,           Jeong-Han%Kim%NULL%0,           Boram%Koh%NULL%2,           Boram%Koh%NULL%0,           Seong Eun%Kim%NULL%2,           Seong Eun%Kim%NULL%0,           Na Ra%Yun%NULL%2,           Na Ra%Yun%NULL%0,           Jae-Hoon%Lee%NULL%2,           Jae-Hoon%Lee%NULL%0,           Jin Yong%Kim%NULL%0,           Jin Yong%Kim%NULL%0,           Yeonjae%Kim%NULL%2,           Yeonjae%Kim%NULL%0,           Ji Hwan%Bang%NULL%2,           Ji Hwan%Bang%NULL%0,           Kyoung-Ho%Song%NULL%3,           Kyoung-Ho%Song%NULL%0,           Hong Bin%Kim%NULL%3,           Hong Bin%Kim%NULL%0,           Ki-hyun%Chung%NULL%2,           Ki-hyun%Chung%NULL%0,           Myoung-don%Oh%NULL%2,           Myoung-don%Oh%NULL%0,           NULL%NULL%NULL%0,           NULL%NULL%NULL%0]'
$ws.Range('E15').Value = '[Wei%Liu%NULL%0,           Zhao-Wu%Tao%NULL%0,           Lei%Wang%NULL%0,           Ming-Li%Yuan%NULL%0,           Kui%Liu%NULL%0,           Ling%Zhou%NULL%0,           Shuang%Wei%NULL%0,           Yan%Deng%NULL%0,           Jing%Liu%NULL%0,           Hui-Guo%Liu%NULL%0,           Ming%Yang%NULL%0,           Yi%Hu%NULL%0,           Pei-Fang%Wei%NULL%0,           Pei-Fang%Wei%NULL%0]'
$ws.Range('E16').Value = '[Yu%Shi%NULL%1,           Xia%Yu%NULL%1,           Hong%Zhao%NULL%1,           Hao%Wang%NULL%2,           Ruihong%Zhao%NULL%0,           Jifang%Sheng%jifang_sheng@zju.edu.cn%0]'
$ws.Range('E17').Value = '[Xiaobo%Yang%NULL%0,           Yuan%Yu%NULL%0,           Jiqian%Xu%NULL%0,           Huaqing%Shu%NULL%0,           Jia''an%Xia%NULL%0,           Hong%Liu%NULL%0,           Yongran%Wu%NULL%0,           Lu%Zhang%NULL%0,           Zhui%Yu%NULL%0,           Minghao%Fang%NULL%0,           Ting%Yu%NULL%0,           Yaxin%Wang%NULL%0,           Shangwen%Pan%NULL%0,           Xiaojing%Zou%NULL%0,           Shiying%Yuan%NULL%0,           You%Shang%NULL%0]'
